$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder two pairs of country names (sharedStrings swap in the source diff) ---
# Row 45/46: Dinamarca <-> Republica Dominicana
$ws.Range("A45").Value = "Republica Dominicana"
$ws.Range("A46").Value = "Dinamarca"

# Row 192/193: Belice <-> Nueva Caledonia
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("A193").Value = "Belice"

# --- Update statistic values (Covid data refresh) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1371176
$ws.Range("C4").Value = 3538
$ws.Range("E4").Value = 1033388
$ws.Range("G4").Value = 91
$ws.Range("H4").Value = 80878

# Row 15: India
$ws.Range("B15").Value = 69294
$ws.Range("C15").Value = 2133
$ws.Range("E15").Value = 45376
$ws.Range("G15").Value = 42
$ws.Range("H15").Value = 2254

# Row 45: Republica Dominicana (new data, after swap)
$ws.Range("B45").Value = 10634
$ws.Range("C45").Value = 287
$ws.Range("D45").Value = 2870
$ws.Range("E45").Value = 7371
$ws.Range("F45").Value = 134
$ws.Range("G45").Value = 5
$ws.Range("H45").Value = 393

# Row 46: Dinamarca (new data, after swap)
$ws.Range("B46").Value = 10513
$ws.Range("C46").Value = 84
$ws.Range("D46").Value = 8328
$ws.Range("E46").Value = 1652
$ws.Range("F46").Value = 43
$ws.Range("G46").Value = 4
$ws.Range("H46").Value = 533

# Row 48: Egipto
$ws.Range("B48").Value = 9746
$ws.Range("C48").Value = 346
$ws.Range("D48").Value = 2172
$ws.Range("E48").Value = 7041
$ws.Range("G48").Value = 8
$ws.Range("H48").Value = 533

# Row 51: Chequia
$ws.Range("B51").Value = 8157
$ws.Range("C51").Value = 34
$ws.Range("D51").Value = 4695
$ws.Range("E51").Value = 3181
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 281

# Row 52: Noruega
$ws.Range("E52").Value = 7866
$ws.Range("G52").Value = 5
$ws.Range("H52").Value = 224

# Row 59: Kazajistan
$ws.Range("D59").Value = 2074
$ws.Range("E59").Value = 3054

# Row 65: Luxemburgo
$ws.Range("B65").Value = 3888
$ws.Range("C65").Value = 2
$ws.Range("D65").Value = 3602
$ws.Range("E65").Value = 185
$ws.Range("F65").Value = 18

# Row 93: Somalia
$ws.Range("B93").Value = 1089
$ws.Range("C93").Value = 35
$ws.Range("D93").Value = 121
$ws.Range("E93").Value = 916
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 52

# Row 148: Birmania
$ws.Range("D148").Value = 74
$ws.Range("E148").Value = 100

# Row 192: Nueva Caledonia (new data, after swap)
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

# Row 193: Belice (new data, after swap)
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2
